$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 90.27778000000001
$ws.Range("I12").Value = 92.57143000000001
$ws.Range("K12").Value = 92.57143000000001
$ws.Range("M12").Value = 77.42856999999999
$ws.Range("H70").Value = 1687
$ws.Range("I70").Value = 1666.3334
$ws.Range("J70").Value = 1749
$ws.Range("K70").Value = 4999.0002
$ws.Range("L70").Value = 5247
$ws.Range("M70").Value = -4729.0002
$ws.Range("N70").Value = -5787
$ws.Range("H73").Value = 1687
$ws.Range("I73").Value = 1666.3334
$ws.Range("J73").Value = 1749
$ws.Range("K73").Value = 4999.0002
$ws.Range("L73").Value = 5247
$ws.Range("M73").Value = -4063.0002
$ws.Range("N73").Value = -7119
$ws.Range("H88").Value = 2262.0908
$ws.Range("J88").Value = 1908.8889
$ws.Range("L88").Value = 1908.8889
$ws.Range("N88").Value = -2720.8889
$ws.Range("H91").Value = 2262.0908
$ws.Range("J91").Value = 1908.8889
$ws.Range("L91").Value = 1908.8889
$ws.Range("N91").Value = -4716.8889
$ws.Range("H98").Value = 1495.7
$ws.Range("I98").Value = 1495.7
$ws.Range("K98").Value = 1495.7
$ws.Range("M98").Value = 2.299999999999955
$ws.Range("H100").Value = 31290.805
$ws.Range("I100").Value = 43075.4
$ws.Range("J100").Value = 4507.636
$ws.Range("K100").Value = 43075.4
$ws.Range("L100").Value = 4507.636
$ws.Range("M100").Value = -42534.4
$ws.Range("N100").Value = -5589.636
$ws.Range("H107").Value = 5627.143
$ws.Range("I107").Value = 1826.8
$ws.Range("J107").Value = 15128
$ws.Range("K107").Value = 1826.8
$ws.Range("L107").Value = 15128
$ws.Range("M107").Value = 93.20000000000005
$ws.Range("N107").Value = -18968
$ws.Range("H122").Value = 1495.7
$ws.Range("I122").Value = 1495.7
$ws.Range("K122").Value = 4487.1
$ws.Range("M122").Value = -2037.1
$ws.Range("H132").Value = 1260.3429
$ws.Range("I132").Value = 1284.875
$ws.Range("J132").Value = 998.6667
$ws.Range("K132").Value = 3854.625
$ws.Range("L132").Value = 2996.0001
$ws.Range("M132").Value = -1324.625
$ws.Range("N132").Value = -8056.0001
$ws.Range("H137").Value = 10663.775
$ws.Range("I137").Value = 3912.6086
$ws.Range("J137").Value = 19797.705
$ws.Range("K137").Value = 11737.8258
$ws.Range("L137").Value = 59393.11500000001
$ws.Range("M137").Value = -9187.825800000001
$ws.Range("N137").Value = -64493.11500000001
$ws.Range("H138").Value = 2240.3948
$ws.Range("I138").Value = 1519.6522
$ws.Range("K138").Value = 4558.9566
$ws.Range("M138").Value = 581.0434000000005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3568.347
$ws.Range("I32").Value = 3132.5
$ws.Range("K32").Value = 3132.5
$ws.Range("M32").Value = -2845.5
$ws.Range("H45").Value = 7654.4346
$ws.Range("I45").Value = 8273.799999999999
$ws.Range("K45").Value = 8273.799999999999
$ws.Range("M45").Value = -7896.799999999999
$ws.Range("H76").Value = 64825.668
$ws.Range("J76").Value = 64825.668
$ws.Range("L76").Value = 64825.668
$ws.Range("N76").Value = -65501.668
$ws.Range("H79").Value = 64825.668
$ws.Range("J79").Value = 64825.668
$ws.Range("L79").Value = 64825.668
$ws.Range("N79").Value = -67165.66800000001
$ws.Range("H102").Value = 1369.5927
$ws.Range("I102").Value = 1399.7142
$ws.Range("J102").Value = 1264.1666
$ws.Range("K102").Value = 1399.7142
$ws.Range("L102").Value = 1264.1666
$ws.Range("M102").Value = 222.2858000000001
$ws.Range("N102").Value = -4508.1666
$ws.Range("H130").Value = 49996
$ws.Range("J130").Value = 49996
$ws.Range("L130").Value = 49996
$ws.Range("N130").Value = -60036
$ws.Range("H132").Value = 1432.2894
$ws.Range("I132").Value = 1327.2812
$ws.Range("J132").Value = 1992.3334
$ws.Range("K132").Value = 3981.8436
$ws.Range("L132").Value = 5977.0002
$ws.Range("M132").Value = -1451.8436
$ws.Range("N132").Value = -11037.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 26059.834
$ws.Range("J103").Value = 26059.834
$ws.Range("L103").Value = 26059.834
$ws.Range("N103").Value = -28403.834
$ws.Range("H134").Value = 7151.4
$ws.Range("I134").Value = 2880.4412
$ws.Range("K134").Value = 8641.3236
$ws.Range("M134").Value = -6106.3236
$ws.Range("H135").Value = 64950
$ws.Range("J135").Value = 64950
$ws.Range("L135").Value = 64950
$ws.Range("N135").Value = -75090
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2096.5908
$ws.Range("I16").Value = 2146.3684
$ws.Range("J16").Value = 1781.3334
$ws.Range("K16").Value = 2146.3684
$ws.Range("L16").Value = 1781.3334
$ws.Range("M16").Value = -1859.3684
$ws.Range("N16").Value = -2355.3334
$ws.Range("H58").Value = 5747.5
$ws.Range("I58").Value = 2098.4
$ws.Range("K58").Value = 2098.4
$ws.Range("M58").Value = -1895.4
$ws.Range("H113").Value = 2096.5908
$ws.Range("I113").Value = 2146.3684
$ws.Range("J113").Value = 1781.3334
$ws.Range("K113").Value = 2146.3684
$ws.Range("L113").Value = 1781.3334
$ws.Range("M113").Value = 23.63160000000016
$ws.Range("N113").Value = -6121.3334
$ws.Range("H132").Value = 19687.328
$ws.Range("I132").Value = 13030.378
$ws.Range("J132").Value = 27897.566
$ws.Range("K132").Value = 39091.13400000001
$ws.Range("L132").Value = 83692.698
$ws.Range("M132").Value = -36561.13400000001
$ws.Range("N132").Value = -88752.698
$ws.Range("H134").Value = 4397.5
$ws.Range("I134").Value = 2447.25
$ws.Range("K134").Value = 7341.75
$ws.Range("M134").Value = -4806.75
$ws.Range("H136").Value = 5747.5
$ws.Range("I136").Value = 2098.4
$ws.Range("K136").Value = 6295.200000000001
$ws.Range("M136").Value = -3745.200000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1150.4546
$ws.Range("J113").Value = 803.25
$ws.Range("L113").Value = 2409.75
$ws.Range("N113").Value = -6749.75
$ws.Range("H132").Value = 4001726.8
$ws.Range("I132").Value = 1363
$ws.Range("K132").Value = 12267
$ws.Range("M132").Value = -9737
$ws.Range("H133").Value = 5999
$ws.Range("J133").Value = 5999
$ws.Range("L133").Value = 17997
$ws.Range("N133").Value = -28117
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1999999
$ws.Range("I113").Value = 1999999
$ws.Range("K113").Value = 1999999
$ws.Range("M113").Value = -1997829
$ws.Range("H132").Value = 13594.974
$ws.Range("I132").Value = 12261.385
$ws.Range("K132").Value = 36784.155
$ws.Range("M132").Value = -34254.155
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2927.7896
$ws.Range("I100").Value = 2268.2222
$ws.Range("J100").Value = 14800
$ws.Range("K100").Value = 2268.2222
$ws.Range("L100").Value = 14800
$ws.Range("M100").Value = -1727.2222
$ws.Range("N100").Value = -15882
$ws.Range("H122").Value = 5068.778
$ws.Range("I122").Value = 6926.25
$ws.Range("J122").Value = 3582.8
$ws.Range("K122").Value = 20778.75
$ws.Range("L122").Value = 10748.4
$ws.Range("M122").Value = -18328.75
$ws.Range("N122").Value = -15648.4
$ws.Range("H132").Value = 10081
$ws.Range("I132").Value = 9997.5
$ws.Range("J132").Value = 10248
$ws.Range("K132").Value = 29992.5
$ws.Range("L132").Value = 30744
$ws.Range("M132").Value = -27462.5
$ws.Range("N132").Value = -35804
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 69949
$ws.Range("J75").Value = 79900
$ws.Range("L75").Value = 79900
$ws.Range("N75").Value = -81772
$ws.Range("H78").Value = 69949
$ws.Range("J78").Value = 79900
$ws.Range("L78").Value = 239700
$ws.Range("N78").Value = -249060
$ws.Range("H87").Value = 39999
$ws.Range("J87").Value = 39999
$ws.Range("L87").Value = 39999
$ws.Range("N87").Value = -42495
$ws.Range("H90").Value = 39999
$ws.Range("J90").Value = 39999
$ws.Range("L90").Value = 119997
$ws.Range("N90").Value = -132477
$ws.Range("H122").Value = 6891.381
$ws.Range("I122").Value = 4428.5
$ws.Range("J122").Value = 14772.6
$ws.Range("K122").Value = 13285.5
$ws.Range("L122").Value = 44317.8
$ws.Range("M122").Value = -10835.5
$ws.Range("N122").Value = -49217.8
$ws.Range("H126").Value = 4302.5757
$ws.Range("I126").Value = 3406.5
$ws.Range("K126").Value = 10219.5
$ws.Range("M126").Value = -7749.5
